$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Power" worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Power"

# Title (bears an underline font style in the original)
$ws2.Range("D14").Value = "Resistor 1"
$ws2.Range("D8").Value = "Resistor 2"
$ws2.Range("F8").Value = "ohm"
$ws2.Range("D6").Value = "Voltage Input"
$ws2.Range("F6").Value = "volts"
$ws2.Range("D7").Value = "Desired Output"
$ws2.Range("F7").Value = "Volts"
$ws2.Range("D10").Value = "Vout / Vin"
$ws2.Range("D11").Value = "1/(R1+R2)"
$ws2.Range("D12").Value = "R1+R2"
$ws2.Range("F15").Value = "amps"
$ws2.Range("F14").Value = "ohms"
$ws2.Range("D15").Value = "Current At Output"
$ws2.Range("D4").Value = "9 volt to 3.3 volt voltage divider"

# Numeric inputs
$ws2.Range("E6").Value = 9
$ws2.Range("E7").Value = 5
$ws2.Range("E8").Value = 10

# Formulas
$ws2.Range("E10").Formula = "=E7/E6"
$ws2.Range("E11").Formula = "=E10/E8"
$ws2.Range("E12").Formula = "=1/E11"
$ws2.Range("E14").Formula = "=E12-E8"
$ws2.Range("E15").Formula = "=E7/E14"

# Style the title with an underline font
$ws2.Range("D4").Font.Underline = $true

# Column width for column D
$ws2.Columns.Item(4).ColumnWidth = 25.88671875

# Sheet1 view changes: scroll so B6 is the top-left visible cell (selection stays at J12)
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 2

# View settings for the new sheet, and make "Power" the active sheet/tab
$ws2.Activate() | Out-Null
$ws2.Range("E8").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130

Write-Output "done"
